$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the two top "number line" diagrams up by one row and left by
#        two columns, which also removes the bold "1)"/"2)" labels that used
#        to sit in column B (this mirrors the author's "cambie las
#        instancias del 2" edit: the old 3rd copy goes away further below). --
$ws.Columns.Item(1).Delete()
$ws.Columns.Item(1).Delete()
$ws.Rows.Item(1).Delete()

# --- 2. Remove everything from row 12 down (the leftover 3rd number-line
#        copy, its gauge swatches, and the old "Curso" legend) so we can
#        rebuild the legend block in its new, higher position. ---
$ws.Range("A12:N31").Clear()

# Restore / set the row heights for the rows below the two diagrams.
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 56.75
$ws.Rows.Item(14).RowHeight = 20
for ($r = 15; $r -le 27; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

# --- 3. Rebuild the "Curso nºX" legend, now starting at row 15. ---
$ws.Range("B15").Interior.Color = 0x0023DC      # RGB(220,35,0) red swatch
$ws.Range("C15").Value = "Curso nº1"
$ws.Range("C15").HorizontalAlignment = -4152     # xlRight

$ws.Range("B17").Interior.Color = 0x47A67D       # RGB(125,166,71) green swatch
$ws.Range("C17").Value = "Curso nº2"
$ws.Range("C17").HorizontalAlignment = -4152

$ws.Range("B19").Interior.Color = 0x0E95FF       # RGB(255,149,14) orange swatch
$ws.Range("C19").Value = "Curso nº3"
$ws.Range("C19").HorizontalAlignment = -4152

$ws.Range("B21").Interior.Color = 0x7B4A60       # purple swatch
$ws.Range("C21").Value = "Curso nº4"
$ws.Range("C21").HorizontalAlignment = -4152

# highlight band to the right of the first two legend rows
$ws.Range("E15:N16").Interior.ColorIndex = 0
$ws.Range("E17:N18").Interior.ColorIndex = 0

$ws.Range("F16").Select()
